# Generate Report for Handback
# Update the Correspond Handoff Datetime / Correspond Handback DateTime
# values for the second handed-back file (9a6476c5-...) on both the
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 3 corresponds to the 9a6476c5-... file
$wsZhCn.Range("E3").Value = "2016-03-20 06:37:04"
$wsZhCn.Range("H3").Value = "2016-03-20 06:37:30"

# de-de sheet: row 3 corresponds to the 9a6476c5-... file
$wsDeDe.Range("E3").Value = "2016-03-20 06:37:07"
$wsDeDe.Range("H3").Value = "2016-03-20 06:37:35"
